$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.477.19"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "1.828.30"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5344"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4018"
$ws.Range("E8").Value = "  +6.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07605"
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.81"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.112"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.340"
$ws.Range("E12").Value = "  +4.33%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.624"
$ws.Range("E13").Value = "  +5.61%  "
$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.9999"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.94"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").Value = "1.828.46"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.56"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001075"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.093"
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").Value = "28.467.41"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.19"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.214"
$ws.Range("E25").Value = "  +5.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.462"
$ws.Range("E26").Value = "  +7.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.55"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.67"
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").Value = "2.039.42"
$ws.Range("E29").Value = "  +2.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.54"
$ws.Range("E30").Value = "  +3.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.124"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07558"
$ws.Range("E33").Value = "  +17.43%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.672"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.641"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2230"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02344"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.210"
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.865"
$ws.Range("E39").Value = "  +5.33%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6268"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.31"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.174"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.392"
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.49"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.703"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5849"
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.85"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.997"
$ws.Range("E49").Value = "  +3.85%  "
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06887"
$ws.Range("E51").Value = "  +1.42%  "
